$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from J1 into K1 so the new header cell matches the bold/
# bordered/centered header formatting used by the rest of row 1.
$ws.Range("J1").Copy($ws.Range("K1"))
$ws.Range("K1").Value = "intervention_type"

$ws.Range("K2").Value = "DRUG"
$ws.Range("K3").Value = "PROCEDURE"
$ws.Range("K4").Value = "OTHER"
$ws.Range("K5").Value = "PROCEDURE"
$ws.Range("K6").Value = "OTHER"
$ws.Range("K7").Value = "RADIATION"
$ws.Range("K8").Value = "BEHAVIORAL"
$ws.Range("K9").Value = "DEVICE"
$ws.Range("K10").Value = "OTHER"
$ws.Range("K11").Value = "DRUG"
$ws.Range("K12").Value = "GENETIC"
$ws.Range("K13").Value = "DRUG"
$ws.Range("K14").Value = "OTHER"
# K15 left blank (matches empty inline string cell in source diff)
$ws.Range("K16").Value = "OTHER"
$ws.Range("K17").Value = "OTHER"
$ws.Range("K18").Value = "BIOLOGICAL"
$ws.Range("K19").Value = "DRUG"
$ws.Range("K20").Value = "BIOLOGICAL"
$ws.Range("K21").Value = "BIOLOGICAL"
$ws.Range("K22").Value = "BIOLOGICAL"
$ws.Range("K23").Value = "PROCEDURE"
$ws.Range("K24").Value = "OTHER"
$ws.Range("K25").Value = "OTHER"
$ws.Range("K26").Value = "OTHER"
$ws.Range("K27").Value = "BIOLOGICAL"
$ws.Range("K28").Value = "BIOLOGICAL"
$ws.Range("K29").Value = "BIOLOGICAL"
$ws.Range("K30").Value = "OTHER"
$ws.Range("K31").Value = "OTHER"
$ws.Range("K32").Value = "OTHER"
$ws.Range("K33").Value = "OTHER"
$ws.Range("K34").Value = "OTHER"
$ws.Range("K35").Value = "BEHAVIORAL"
$ws.Range("K36").Value = "OTHER"
$ws.Range("K37").Value = "OTHER"
$ws.Range("K38").Value = "OTHER"
$ws.Range("K39").Value = "GENETIC"
$ws.Range("K40").Value = "OTHER"
$ws.Range("K41").Value = "PROCEDURE"
$ws.Range("K42").Value = "BEHAVIORAL"
$ws.Range("K43").Value = "BIOLOGICAL"
$ws.Range("K44").Value = "DIAGNOSTIC_TEST"
$ws.Range("K45").Value = "DEVICE"
$ws.Range("K46").Value = "OTHER"
$ws.Range("K47").Value = "OTHER"
$ws.Range("K48").Value = "DRUG"
$ws.Range("K49").Value = "OTHER"
$ws.Range("K50").Value = "OTHER"
$ws.Range("K51").Value = "BIOLOGICAL"
$ws.Range("K52").Value = "RADIATION"
$ws.Range("K53").Value = "OTHER"
$ws.Range("K54").Value = "PROCEDURE"
$ws.Range("K55").Value = "DRUG"
$ws.Range("K56").Value = "DRUG"
$ws.Range("K57").Value = "OTHER"
$ws.Range("K58").Value = "OTHER"
$ws.Range("K59").Value = "BEHAVIORAL"
$ws.Range("K60").Value = "OTHER"
$ws.Range("K61").Value = "DRUG"
$ws.Range("K62").Value = "OTHER"
$ws.Range("K63").Value = "DEVICE"
$ws.Range("K64").Value = "OTHER"
$ws.Range("K65").Value = "PROCEDURE"
